function Set-TextValue {
    param($Cell, $Text)
    if ($Text -match '^-?\d+(\.\d+)?$') {
        $Cell.Value = "'" + $Text
    } else {
        $Cell.Value = $Text
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-level Price (D) / Volume(1h) (E) updates
$rowUpdates = @(
  @{ Row=2; D="62.775.09"; E="  +1.12%  " },
  @{ Row=3; D="3.471.62"; E="  +1.25%  " },
  @{ Row=4; E="  -0.28%  " },
  @{ Row=5; D="414.09"; E="  +1.21%  " },
  @{ Row=6; D="130.81"; E="  +0.49%  " },
  @{ Row=7; E="  -0.98%  " },
  @{ Row=8; E="  +0.04%  " },
  @{ Row=9; D="0.728"; E="  -1.77%  " },
  @{ Row=10; E="  +4.31%  " },
  @{ Row=11; D="42.70"; E="  -0.36%  " },
  @{ Row=12; D="9.61"; E="  +3.80%  " },
  @{ Row=13; E="  -1.41%  " },
  @{ Row=14; D="4.022.39"; E="  +1.26%  " },
  @{ Row=15; E="  -0.22%  " },
  @{ Row=16; D="20.55"; E="  -3.72%  " },
  @{ Row=17; D="3.447.88"; E="  +0.91%  " },
  @{ Row=18; D="12.64"; E="  +0.72%  " },
  @{ Row=19; E="  -1.43%  " },
  @{ Row=20; D="62.755.86"; E="  +1.17%  " },
  @{ Row=21; D="464.37"; E="  +1.74%  " },
  @{ Row=22; D="90.82"; E="  -0.89%  " },
  @{ Row=23; E="  +1.50%  " },
  @{ Row=24; D="13.28"; E="  +1.30%  " },
  @{ Row=25; D="10.67"; E="  +17.21%  " },
  @{ Row=26; D="3.32"; E="  +1.80%  " },
  @{ Row=27; D="33.40" },
  @{ Row=28; E="  +0.48%  " },
  @{ Row=29; D="7.57"; E="  -2.03%  " },
  @{ Row=30; E="  -0.76%  " },
  @{ Row=31; D="2.66"; E="  -1.08%  " },
  @{ Row=32; D="0.166"; E="  -2.56%  " },
  @{ Row=33; E="  -1.62%  " },
  @{ Row=34; D="40.91"; E="  -4.98%  " },
  @{ Row=35; E="  +0.06%  " },
  @{ Row=36; D="58.71"; E="  +7.85%  " },
  @{ Row=37; E="  -2.60%  " },
  @{ Row=40; D="147.41"; E="  +3.37%  " },
  @{ Row=41; E="  -0.56%  " },
  @{ Row=42; E="  +0.32%  " },
  @{ Row=43; E="  -1.47%  " },
  @{ Row=44; E="  +7.19%  " },
  @{ Row=45; E="  +3.74%  " },
  @{ Row=46; D="4.35"; E="  +2.15%  " },
  @{ Row=47; D="2.42"; E="  +13.07%  " },
  @{ Row=48; D="0.0₃0552"; E="  +28.21%  " },
  @{ Row=49; D="16.40"; E="  -1.67%  " },
  @{ Row=50; D="22.35"; E="  -0.06%  " },
  @{ Row=51; D="0.141"; E="  +1.14%  " }
)

foreach ($r in $rowUpdates) {
    if ($r.ContainsKey("D")) {
        Set-TextValue $ws.Cells.Item($r.Row, 4) $r.D
    }
    if ($r.ContainsKey("E")) {
        Set-TextValue $ws.Cells.Item($r.Row, 5) $r.E
    }
}

# Rows 38 and 39: Stacks and FirstDigitalUSD swapped rank position, with updated price/volume
Set-TextValue $ws.Cells.Item(38, 2) "Stacks"
Set-TextValue $ws.Cells.Item(38, 3) "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Cells.Item(38, 4) "3.08"
Set-TextValue $ws.Cells.Item(38, 5) "  +4.40%  "

Set-TextValue $ws.Cells.Item(39, 2) "FirstDigitalUSD"
Set-TextValue $ws.Cells.Item(39, 3) "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Cells.Item(39, 4) "1.00"
Set-TextValue $ws.Cells.Item(39, 5) "  +0.15%  "
